# cryptos.xlsx refresh - GitHub Actions data pull (prices / 1h volume).
# Source feed re-scraped; most rows get new "Price" / "Volume(1h)" text,
# two rows (44/45) swap rank between OKB and VeChain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these columns hold plain text (e.g. "140.10", "  +1.92%  ") even
# though some look numeric. Assigning such a string straight to .Value makes
# Excel auto-convert it to a real number (losing the trailing zero / exact
# text), so for the handful of digit-only replacements we instead build it
# as a text formula and paste-special the computed value back over itself -
# that yields a genuine text cell without leaving a "Text" number-format/
# quote-prefix mark on the cell (keeping cell formatting untouched, as in
# the source edit).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $escaped = $val -replace '"', '""'
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163) # xlPasteValues
}

$ws.Range('D2').Value = '58.724.68'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').Value = '3.153.22'
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '534.47'
$ws.Range('E5').Value = '  +1.21%  '
Set-TextValue 'D6' '140.10'
$ws.Range('E6').Value = '  +1.92%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  +10.65%  '
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('E10').Value = '  +3.00%  '
$ws.Range('E11').Value = '  +4.30%  '
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').Value = '3.696.41'
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('E14').Value = '  +1.81%  '
$ws.Range('E15').Value = '  +5.90%  '
$ws.Range('D16').Value = '58.755.36'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').Value = '3.153.35'
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('E18').Value = '  +4.90%  '
Set-TextValue 'D19' '13.00'
$ws.Range('E19').Value = '  +3.79%  '
$ws.Range('E20').Value = '  +3.63%  '
Set-TextValue 'D21' '372.44'
$ws.Range('E21').Value = '  +6.48%  '
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('E23').Value = '  +0.04%  '
Set-TextValue 'D24' '69.69'
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  +0.21%  '
Set-TextValue 'D28' '8.00'
$ws.Range('E28').Value = '  +12.16%  '
$ws.Range('E29').Value = '  +1.23%  '
Set-TextValue 'D30' '6.17'
$ws.Range('E30').Value = '  +2.94%  '
$ws.Range('E31').Value = '  +1.51%  '
Set-TextValue 'D32' '21.91'
$ws.Range('E32').Value = '  +3.71%  '
Set-TextValue 'D33' '5.19'
$ws.Range('E33').Value = '  +6.83%  '
$ws.Range('E34').Value = '  +2.70%  '
Set-TextValue 'D35' '160.02'
$ws.Range('E35').Value = '  +0.61%  '
Set-TextValue 'D36' '6.26'
$ws.Range('E36').Value = '  +4.00%  '
$ws.Range('E37').Value = '  +9.92%  '
Set-TextValue 'D38' '25.31'
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').Value = '2.655.19'
$ws.Range('E39').Value = '  +11.05%  '
Set-TextValue 'D40' '1.68'
$ws.Range('E40').Value = '  +4.72%  '
$ws.Range('E41').Value = '  +3.14%  '
Set-TextValue 'D42' '4.18'
$ws.Range('E42').Value = '  +3.90%  '
$ws.Range('E43').Value = '  +2.50%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D44' '0.0286'
$ws.Range('E44').Value = '  +7.81%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D45' '38.53'
$ws.Range('E45').Value = '  +4.40%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '3.194.53'
$ws.Range('E47').Value = '  +1.98%  '
$ws.Range('E48').Value = '  +13.61%  '
Set-TextValue 'D49' '0.984'
$ws.Range('E49').Value = '  +3.35%  '
$ws.Range('E50').Value = '  +2.87%  '
Set-TextValue 'D51' '20.25'
$ws.Range('E51').Value = '  +3.84%  '

$excel.CutCopyMode = $false
